# Rolling the forecast window forward by one week:
#  - "Forecast Comparison": Week_Start_Date (col B) shifts to the next week,
#    and MyForecast (col D) gets refreshed values.
#  - "Summary": recomputed stats reflecting the new forecast window.
#
# Helper: write a value as TEXT (preserve the original inline/shared-string
# cell type instead of letting Excel auto-coerce numeric- or date-looking
# strings into real numbers/dates), then drop back to the Normal style so no
# stray number-format style lingers on the cell.
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Forecast Comparison
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$forecastRows = @(
    @{ Row = 2;  Date = "2025-01-12"; MyForecast = 749 },
    @{ Row = 3;  Date = "2025-01-19"; MyForecast = 739 },
    @{ Row = 4;  Date = "2025-01-26"; MyForecast = 742 },
    @{ Row = 5;  Date = "2025-02-02"; MyForecast = 775 },
    @{ Row = 6;  Date = "2025-02-09"; MyForecast = 517 },
    @{ Row = 7;  Date = "2025-02-16"; MyForecast = 508 },
    @{ Row = 8;  Date = "2025-02-23"; MyForecast = 470 },
    @{ Row = 9;  Date = "2025-03-02"; MyForecast = 699 },
    @{ Row = 10; Date = "2025-03-09"; MyForecast = 583 },
    @{ Row = 11; Date = "2025-03-16"; MyForecast = 529 },
    @{ Row = 12; Date = "2025-03-23"; MyForecast = 533 },
    @{ Row = 13; Date = "2025-03-30"; MyForecast = 532 },
    @{ Row = 14; Date = "2025-04-06"; MyForecast = 495 },
    @{ Row = 15; Date = "2025-04-13"; MyForecast = 453 },
    @{ Row = 16; Date = "2025-04-20"; MyForecast = 470 },
    @{ Row = 17; Date = "2025-04-27"; MyForecast = 536 }
)

foreach ($entry in $forecastRows) {
    $dateCell = $wsForecast.Cells.Item($entry.Row, 2)   # column B: Week_Start_Date
    Set-TextValue $dateCell $entry.Date

    $forecastCell = $wsForecast.Cells.Item($entry.Row, 4) # column D: MyForecast
    $forecastCell.Value = $entry.MyForecast
}

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

Set-TextValue $wsSummary.Range("B2")  "2023-01-01 to 2025-01-05"   # Historical Range
Set-TextValue $wsSummary.Range("B4")  "1033"                        # Max Sales
Set-TextValue $wsSummary.Range("B5")  "399"                         # Mean Sales
Set-TextValue $wsSummary.Range("B6")  "363"                         # Median Sales
Set-TextValue $wsSummary.Range("B7")  "244"                         # Std Dev Sales
Set-TextValue $wsSummary.Range("B8")  "41117 units"                 # Total Historical Sales
Set-TextValue $wsSummary.Range("B9")  "9330"                        # Total Forecast (16 Weeks)
Set-TextValue $wsSummary.Range("B10") "5199"                        # Total Forecast (8 Weeks)
Set-TextValue $wsSummary.Range("B11") "3005"                        # Total Forecast (4 Weeks)
Set-TextValue $wsSummary.Range("B12") "775"                         # Max Forecast
Set-TextValue $wsSummary.Range("B13") "2025-02-02"                  # Max Forecast Week
Set-TextValue $wsSummary.Range("B14") "453"                         # Min Forecast
# Min Forecast Week (B15) stays "2025-04-13" - unchanged by this edit.
